$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Force the new "2020" label into the shared-string table as text (it would
# otherwise be auto-coerced to a number since it looks numeric), matching
# how the other year labels (A2:A21) are stored, then drop the temporary
# number-format style so the cell keeps the sheet's default style.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "2020"
$ws.Range("A22").Style = "Normal"

$ws.Range("B22").Value = 30.7
